$d = $word.ActiveDocument

$replacements = @(
    @("59×18=", "87×20="),
    @("81×47=", "54×76="),
    @("37×35=", "64×69="),
    @("38×82=", "82×60="),
    @("38×31=", "73×16="),
    @("83×44=", "85×70="),
    @("87×14=", "88×92="),
    @("12×22=", "38×59="),
    @("84×16=", "29×12="),
    @("98×96=", "82×15="),
    @("71×55=", "52×44="),
    @("30×29=", "27×29="),
    @("65×49=", "24×63="),
    @("76×27=", "86×33="),
    @("55×65=", "83×88="),
    @("72×85=", "56×35="),
    @("15×91=", "28×86="),
    @("80×98=", "87×93="),
    @("48×23=", "30×93="),
    @("29×28=", "46×53="),
    @("68×99=", "33×81="),
    @("30×82=", "76×99="),
    @("26×54=", "66×79="),
    @("47×45=", "47×16="),
    @("81×82=", "19×56=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
